# This workbook holds a weekly time series of "Poroto granado" price
# records at the top of the data block (rows 2..139), ordered with the
# most recent week first. A new weekly record is inserted at row 33,
# pushing the existing rows 33..139 down to rows 34..140.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at row 33 (shifts rows 33:139 down to 34:140).
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new weekly record.
$ws.Range("A33").Value = 8
$ws.Range("B33").Value = "Terminal La Palmera de La Serena"
$ws.Range("C33").Value = "Coquimbo"
$ws.Range("D33").Value = 45054
$ws.Range("E33").Value = 4
$ws.Range("F33").Value = 100112030
$ws.Range("G33").Value = "Poroto granado"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 480
$ws.Range("K33").Value = 29000
$ws.Range("L33").Value = 30000
$ws.Range("M33").Value = 29500
$ws.Range("N33").Value = "`$/malla 25 kilos"
$ws.Range("O33").Value = "Provincia de Limarí"
$ws.Range("P33").Value = 1180
$ws.Range("Q33").Value = 25
$ws.Range("R33").Value = "Hortaliza"
